# Update output data and maps: replace B2:D9 correlation values with new computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @(-0.8907367373884217, 0.7669017931388056, 0.7428693640827839)
    3 = @(0.5962233087516597, -0.5704045726661989, 0.7137775600844258)
    4 = @(-0.7461621130731895, 0.6410763752544432, 0.6548595919230661)
    5 = @(0.689057761770427, 0.6145248210607654, 0.4847682382022319)
    6 = @(0.635789432315277, -0.6901269622900971, 0.5900641550764121)
    7 = @(0.7542379233426603, -0.5802153095351402, 0.7239421788580264)
    8 = @(-0.7145350778500471, 0.7747217287508602, -0.6822174105333796)
    9 = @(0.7885243520054063, -0.7807334039691967, -0.6944669369326431)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
}
